# Updating xpaths and adding wait for buttons
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column A (Button) and set a width for the new column C (SelectValue)
# so the longer xpath/selector strings are easier to read.
$ws.Columns.Item(1).ColumnWidth = 23.1666666666667
$ws.Columns.Item(3).ColumnWidth = 41.3072916666667

# Move the active selection to G8 (single cell) after reviewing the data.
$ws.Range("G8").Select() | Out-Null
